# Update the Ccl12-Ccr3 LR-pairs worksheet with new TPM-derived values.
# A new "ECs" sending cluster is introduced (2 new rows at the top of the
# data block), and every existing data row's computed metrics are
# refreshed to match the new TPM run. Row 1 (the header) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Ccl12", "Ccr3", "Inflammatory-Mac", 1, 0.3333333333333333, 0.2401406666666667, 0.720422, 0.002192729333974893, 0.002203844280310985, 1, 0.3333333333333333, 0.1375686666666667, 0.412706, 0.2896572731203081, 0.2896572731203081, 0.03303583132577778, 0.297322481932, 0.0006351399995700769, 0.0006383595246166681),
    @("ECs", "Ccl12", "Ccr3", "Resolving-Mac", 1, 0.3333333333333333, 0.2401406666666667, 0.720422, 0.002192729333974893, 0.002203844280310985, 2, 0.6666666666666666, 0.3373673333333334, 1.012102, 0.7103427268796919, 0.7103427268796919, 0.08101561633822223, 0.729140547044, 0.001557589334404816, 0.001565484755694317),
    @("FAPs", "Ccl12", "Ccr3", "Inflammatory-Mac", 1, 0.3333333333333333, 0.643692, 1.931076, 0.005877564804149375, 0.005907358183739276, 1, 0.3333333333333333, 0.1375686666666667, 0.412706, 0.2896572731203081, 0.2896572731203081, 0.08855185018400001, 0.796966651656, 0.001702479393757806, 0.001711109262846855),
    @("FAPs", "Ccl12", "Ccr3", "Resolving-Mac", 1, 0.3333333333333333, 0.643692, 1.931076, 0.005877564804149375, 0.005907358183739276, 2, 0.6666666666666666, 0.3373673333333334, 1.012102, 0.7103427268796919, 0.7103427268796919, 0.217160653528, 1.954445881752, 0.004175085410391569, 0.004196248920892422),
    @("Inflammatory-Mac", "Ccl12", "Ccr3", "Inflammatory-Mac", 3, 1, 60.66100033333333, 181.983001, 0.553896833491318, 0.5567045368793268, 1, 0.3333333333333333, 0.1375686666666667, 0.412706, 0.2896572731203081, 0.2896572731203081, 8.34505293452289, 75.105476410706, 0.1604402463790685, 0.1612535180861698),
    @("Inflammatory-Mac", "Ccl12", "Ccr3", "Resolving-Mac", 3, 1, 60.66100033333333, 181.983001, 0.553896833491318, 0.5567045368793268, 2, 0.6666666666666666, 0.3373673333333334, 1.012102, 0.7103427268796919, 0.7103427268796919, 20.46503991978911, 184.185359278102, 0.3934565871122495, 0.395451018793157),
    @("MuSCs", "Ccl12", "Ccr3", "Inflammatory-Mac", 1, 0.5, 1.6570225, 3.314045, 0.0151303063043872, 0.01013800122420362, 1, 0.3333333333333333, 0.1375686666666667, 0.412706, 0.2896572731203081, 0.2896572731203081, 0.2279543759616667, 1.36772625577, 0.004382603265603804, 0.002936545789493166),
    @("MuSCs", "Ccl12", "Ccr3", "Resolving-Mac", 1, 0.5, 1.6570225, 3.314045, 0.0151303063043872, 0.01013800122420362, 2, 0.6666666666666666, 0.3373673333333334, 1.012102, 0.7103427268796919, 0.7103427268796919, 0.5590252620983334, 3.35415157259, 0.0107477030387834, 0.007201455434710454),
    @("Resolving-Mac", "Ccl12", "Ccr3", "Inflammatory-Mac", 3, 1, 46.31492933333334, 138.944788, 0.4229025660661706, 0.4250462594324195, 1, 0.3333333333333333, 0.1375686666666667, 0.412706, 0.2896572731203081, 0.2896572731203081, 6.371483075147556, 57.34334767632801, 0.1224968040823079, 0.1231177404571817),
    @("Resolving-Mac", "Ccl12", "Ccr3", "Resolving-Mac", 3, 1, 46.31492933333334, 138.944788, 0.4229025660661706, 0.4250462594324195, 2, 0.6666666666666666, 0.3373673333333334, 1.012102, 0.7103427268796919, 0.7103427268796919, 15.62514420270845, 140.626297824376, 0.3004057619838627, 0.3019285189752378),
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
